$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear previous contents (wipes stale per-row heights from the old layout too) ---
$ws.Rows("1:7").Delete()

# --- Column sizing ---
$ws.Columns("A").ColumnWidth = 16.3

# --- Data ---
$data = @(
  @("Sonde 3.1",    "GTTCCGTTTGACAATAAAAAGGGATATGCG\colorbox{Snow2}{12}A\colorbox{Snow2}{3}TGTATTGTCCYTGAGAG", "559 bp"),
  @("Sonde 1.1",    "GCTAAACTTGTTGCTACTGATGATCTTACAG\colorbox{Snow2}{12}G\colorbox{Snow2}{3}AGGATGAAGAAGATGG", "707 bp"),
  @("forward 1",    "TTGAATGCATATGACCAGAGTGGAAGGCTT", "654 bp"),
  @("forward 3",    "GGAGGAAGTAAACACTCAGAAAGAAGGGAA", "509 bp"),
  @("forward 4",    "GAGACATGAACAACAGAGATGCAAGGCAAA", "472 bp"),
  @("reverse 1.1",  "CCTTCATTAAGACGCTCGAAGAGTGARTTG", "720 bp"),
  @("reverse 1.2",  "TGAATGTCCTTCATTAAGACGCTCGAAGAG", "727 bp"),
  @("reverse 3.3",  "ATTGGGGTGTTTGAGGAATGTTCCGTTTAC", "565 bp"),
  @("reverse 3.6",  "ATCCATTGGGGTGTTTGAGGAATGTTCCGT", "569 bp"),
  @("reverse 3.10", "CTTGTATCCATTGGGGTGTTTGAGGAATGT", "574 bp"),
  @("reverse 3.15", "AGTTGATAAGGACTTGTATCCATTGGGGTG", "586 bp")
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $i + 1
  $ws.Cells.Item($row, 1).Value = $data[$i][0]
  $ws.Cells.Item($row, 2).Value = $data[$i][1]
  $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# --- Formatting: column A (label column) ---
# (engine's Range() doesn't support comma multi-area unions, so apply per block)
foreach ($addr in @("A1:A11", "A14:A15")) {
  $colA = $ws.Range($addr)
  $colA.Font.Name = "Arial"
  $colA.Font.Size = 10
  $colA.Interior.Pattern = 1
  $colA.Interior.ThemeColor = 2
  $colA.Interior.TintAndShade = 0
  $colA.Borders.LineStyle = 1
  $colA.Borders.Weight = 2
  $colA.HorizontalAlignment = -4131
  $colA.VerticalAlignment = -4108
}

# --- Formatting: column B (sequence column) ---
foreach ($addr in @("B1:B11", "B14:B15")) {
  $colB = $ws.Range($addr)
  $colB.Font.Name = "Arial"
  $colB.Font.Size = 10
  $colB.Interior.Pattern = 1
  $colB.Interior.ThemeColor = 2
  $colB.Interior.TintAndShade = 0
  $colB.Borders.LineStyle = 1
  $colB.Borders.Weight = 2
  $colB.VerticalAlignment = -4108
}

# --- View state ---
$excel.ActiveWindow.Zoom = 66
$ws.Range("F17").Select()
